$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as plain text (matching the
# original sheet, where every Price cell is stored as a text string).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = "26.718.03"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.599.27"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "211.43"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("D12").Value = "1.824.22"

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.04"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.555.02"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "65.34"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "26.690.48"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").Value = "0.0₃0759"
$ws.Range("E18").Value = "  +4.08%  "

$ws.Range("D19").Value = "209.78"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").Value = "7.15"
$ws.Range("E21").Value = "  +3.97%  "

$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("D25").Value = "143.11"
$ws.Range("E25").Value = "  -1.81%  "

$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").Value = "0.0519"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").Value = "1.289.00"
$ws.Range("E34").Value = "  +0.55%  "

$ws.Range("E35").Value = "  -5.38%  "

$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").Value = "0.0171"
$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("E39").Value = "  +17.38%  "

$ws.Range("D40").Value = "0.827"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").Value = "5.44"
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.784"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "63.08"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("D45").Value = "1.736.78"
$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").Value = "91.32"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("D47").Value = "1.57"
$ws.Range("E47").Value = "  -1.17%  "

$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  -0.52%  "

$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("E51").Value = "  +0.11%  "

